{"js": "// Word Online (Office.js) edit script.\n// Applies two textual changes to the \"Example questions\" / \"Longer MCQs\" sections:\n//   1. \"The square root of 9 is: _\"  ->  \"The square root of 16 is: _\"\n//   2. A 3-way rotation of the answer options under\n//      \"What is true about a 95% confidence interval of the mean?\":\n//        \"there is a 95% probability that the true mean lies within this range\"\n//          -> \"if you repeated the process many times, 95% of intervals calculated in this way contain the true mean\"\n//        \"if you repeated the process many times, 95% of intervals calculated in this way contain the true mean\"\n//          -> \"95% of the data fall within this range\"\n//        \"95% of the data fall within this range\"\n//          -> \"there is a 95% probability that the true mean lies within this range\"\n//\n// Each of these four strings is unique in the document, so `body.search`\n// can locate the exact (single) run of text to rewrite. All four old\n// strings are searched for *before* any edit is written, and the\n// rotation's new values are computed purely from the old text -- so a\n// replacement written by one step can never be re-discovered and\n// clobbered by a later search/replace (the classic hazard with chained\n// find-and-replace on a 3-way rotation where each new value contains the\n// \"95%\" marker used to find the next target).\n\nconst body = context.document.body;\n\nconst OLD_SQRT = \"The square root of 9 is: _\";\nconst NEW_SQRT = \"The square root of 16 is: _\";\n\nconst TEXT_A = \"there is a 95% probability that the true mean lies within this range\";\nconst TEXT_B = \"if you repeated the process many times, 95% of intervals calculated in this way contain the true mean\";\nconst TEXT_C = \"95% of the data fall within this range\";\n\n// Rotation computed from the *original* text: A -> B, B -> C, C -> A.\nconst replacements = [\n  [OLD_SQRT, NEW_SQRT],\n  [TEXT_A, TEXT_B],\n  [TEXT_B, TEXT_C],\n  [TEXT_C, TEXT_A],\n];\n\n// Locate every old string first (all searches are against the\n// still-unmodified document).\nconst searchResults = replacements.map(([oldText]) =>\n  body.search(oldText, { matchCase: true, matchWholeWord: false })\n);\nsearchResults.forEach((r) => r.load(\"items\"));\nawait context.sync();\n\n// Only now perform the writes, using the ranges found above.\nfor (let i = 0; i < replacements.length; i++) {\n  const [, newText] = replacements[i];\n  const items = searchResults[i].items;\n  for (let j = 0; j < items.length; j++) {\n    items[j].insertText(newText, \"Replace\");\n  }\n}\n\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n# Applies two textual changes to the \"Example questions\" / \"Longer MCQs\" sections:\n#   1. \"The square root of 9 is: _\"  ->  \"The square root of 16 is: _\"\n#   2. A 3-way rotation of the answer options under\n#      \"What is true about a 95% confidence interval of the mean?\":\n#        \"there is a 95% probability that the true mean lies within this range\"\n#          -> \"if you repeated the process many times, 95% of intervals calculated in this way contain the true mean\"\n#        \"if you repeated the process many times, 95% of intervals calculated in this way contain the true mean\"\n#          -> \"95% of the data fall within this range\"\n#        \"95% of the data fall within this range\"\n#          -> \"there is a 95% probability that the true mean lies within this range\"\n#\n# The three \"95%\" strings are rotated simultaneously (based on each\n# paragraph's *original* text) rather than via sequential Find/Replace,\n# so that a replacement written by one step can never be re-matched and\n# clobbered by a later step.\n\n$d = $word.ActiveDocument\n\n$OLD_SQRT = \"The square root of 9 is: _\"\n$NEW_SQRT = \"The square root of 16 is: _\"\n\n$TEXT_A = \"there is a 95% probability that the true mean lies within this range\"\n$TEXT_B = \"if you repeated the process many times, 95% of intervals calculated in this way contain the true mean\"\n$TEXT_C = \"95% of the data fall within this range\"\n\n$rotation = @{}\n$rotation[$TEXT_A] = $TEXT_B\n$rotation[$TEXT_B] = $TEXT_C\n$rotation[$TEXT_C] = $TEXT_A\n\n$count = $d.Paragraphs.Count\nfor ($i = 1; $i -le $count; $i++) {\n    $para = $d.Paragraphs.Item($i)\n    $raw = $para.Range.Text\n    # Paragraph.Range.Text always ends with the paragraph mark (CR); it may\n    # also carry a manual line break (VT) just before it. Strip those off\n    # for comparison, then write back only the visible text so the mark(s)\n    # are left untouched.\n    $clean = $raw.TrimEnd([char]13, [char]11)\n\n    if ($clean -eq $OLD_SQRT) {\n        $para.Range.Text = $NEW_SQRT\n    } elseif ($rotation.ContainsKey($clean)) {\n        $para.Range.Text = $rotation[$clean]\n    }\n}\n"}
